# Rename worksheet 1 ("Property1") to "DataNode" -- part of the commit's
# conceptual unification of DataNode / DataTable / Entity naming.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "DataNode"

# Move the active cell/selection to W37 (matches the author's saved
# selection state in the commit).
$ws.Activate()
$ws.Range("W37").Select()
